# This workbook contains a weekly price log where each new week's record is
# inserted at the top of the data block (row 17, just below the fixed header
# rows 1-16) and all previously existing data rows shift down by one.
#
# Insert a new row at row 17. This pushes old rows 17-45 down to rows 18-46
# and automatically carries over the date-formatted style (s="2") already
# present on column D from the row above, matching the original file's
# per-cell styling instead of a full-row style/format copy.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly record.
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = "Femacal de La Calera"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44497
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100108
$ws.Range("H17").Value = "Tropicales y subtropicales"
$ws.Range("I17").Value = 100108004
$ws.Range("J17").Value = "Papaya"
$ws.Range("K17").Value = "Cultivar IV Región"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 67
$ws.Range("N17").Value = 22000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 22000
$ws.Range("Q17").Value = "`$/bandeja 10 kilos"
$ws.Range("R17").Value = "Provincia del Elquí"
$ws.Range("S17").Value = 2200
$ws.Range("T17").Value = 10
